$d = $word.ActiveDocument

# --- Change 1: "Programa" paragraph ---
# "Programa1. TÉCNICAS..." -> one run containing "Programa" + two manual line
# breaks + "1. TÉCNICAS...". Locate the whole paragraph by its unique text,
# clear it, then splice in the exact target run/break structure as raw OOXML.
$r1 = $d.Content
$r1.Find.Execute("Programa1. TÉCNICAS", $false, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$p1 = $r1.Paragraphs(1).Range
[void]$p1.MoveEnd(1, -1)
$p1.Text = ""
[void]$p1.InsertXML('<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:r><w:t>Programa</w:t><w:br/><w:br/><w:t>1. TÉCNICAS DE JUNÇÃO DE MATERIAIS: Razões técnicas para a junção de materiais, junção por difusão, brasagem, soldagem por explosão, elementos de fixação. 2. PROCESSOS DE SOLDAGEM: Definição de soldagem por fusão, física da soldagem, principais processos de soldagem. 3. NOMENCLATURA DAS JUNTAS SOLDADAS: Desenho e simbologia para soldagem, símbolos básicos, tipos de juntas e soldas, simbologia para soldas em desenho. 4. METALURGIA FÍSICA DAS REGIÕES SOLDADAS: metalurgia da soldagem, estruturas de solidificação, transformações de fase pós-soldagem, transformações de fases em juntas de aço soldadas, ligas de alumínio, ligas de cobre e em metais e ligas especiais. 5. SEGURANÇA NO PROCESSO DE SOLDAGEM: Problemas associados à vaporização de metais, luminosidade, calor e eletricidade. 6. APLICAÇÕES DE JUNTAS SOLDADAS EM ENGENHARIA: Exemplos de estruturas soldadas em engenharia, descontinuidades e defeitos de soldagem, métodos de inspeção em soldas. 7. PRÁTICA EXPERIMENTAL SUPERVISIONADA: Caracterização microestrutural de juntas soldadas (materiais e processos a serem definidos na ocasião da prática experimental), incluindo a redação de relatório técnico de cada grupo. 8. Pós Metálicos - obtenção por processos químicos, termoquímicos, eletrolíticos, atomização e moagem, Caracterização de pós e sua aplicação na metalurgia do pó. 9. Técnicas de mistura, aspectos sobre o transporte e armazenamento de pós, 10. Processos de fabricação de peças "verdes" por compactação uniaxial e isostática, 11. Técnicas de sinterização e fenômenos envolvidos, 12. Manufatura aditiva (impressão 3D). 13. Fontes de calor (laser e feixe de elétrons: obtenção e aplicações), características desejáveis dos pós, parâmetros relevantes no processamento a laser e por feixe de elétrons. 14. Estudos de casos e comparação entre as técnicas estudadas no semestre.</w:t></w:r></w:p>')

# --- Change 2: "Bibliografia" paragraph ---
# Split the single run of numbered references into one run containing 8
# <w:t> segments separated by manual line breaks (<w:br/>), preserving the
# exact original wording/whitespace of each segment.
$r2 = $d.Content
$r2.Find.Execute("1. WAINER, E. et al.", $false, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$p2 = $r2.Paragraphs(1).Range
[void]$p2.MoveEnd(1, -1)
$p2.Text = ""
[void]$p2.InsertXML('<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:r><w:t xml:space="preserve">1. WAINER, E. et al. Soldagem - Processos e Metalurgia. São Paulo: Edgar Blücher, 1992. 494 p. </w:t><w:br/><w:t xml:space="preserve">2. QUITES, A. M., DUTRA, J. C. Tecnologia da soldagem a arco voltaico. Florianópolis: EDEME, 1979. 248 p. </w:t><w:br/><w:t xml:space="preserve">3. GOURD, L. M. Principles of welding technology. London: Edward Arnold, 1980. 218 p. </w:t><w:br/><w:t xml:space="preserve">4. KOU, S. Welding metallurgy, 2nd ed.: John Wiley &amp; Sons, 2003. 461 p. </w:t><w:br/><w:t>5. MESSLER, Jr. R. W. Principles of welding: Processes, physics, chemistry and metallurgy: Wiley VCH Verlag GmbH &amp; Co., 2004. 662 p.</w:t><w:br/><w:t>6. KALPAKJIAN, S.; SCHMID, S. Manufacturing processes for engineering materials. 5ª ed., Pearson Education, New Jersey, 2007.</w:t><w:br/><w:t>7. GERMAN, R.M. Sintering theory and practice. New York, Wiley-Interscience, 1996</w:t><w:br/><w:t>8. GIBSON, I., ROSEN, D., STUCKER, B., Additive Manufacturing Technologies, New York, Springer Verlag, 2015.</w:t></w:r></w:p>')
